$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'Dr. Gehan Adel, Dr. Veronia Rafat, Dr. Amira Sobhy, Administrator, Dr. Servinaz Sayed Mohammad'
$ws.Range("G3").Value = 'Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Majorelle Magdy, Administrator'
$ws.Range("G4").Value = 'Dr. Gehan Adel, Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Majorelle Magdy, Dr. Servinaz Sayed Mohammad'
$ws.Range("G5").Value = 'Dr. Amira Sobhy, Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Veronia Rafat'
$ws.Range("G6").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Mohammad El-Tanany, Dr. Manar Montaser, Dr. Majorelle Magdy, Dr. Alshimaa Atef'
$ws.Range("G7").Value = 'Dr. Nada Mohammad, Dr. Amera Ahmad Saad, Dr. Lamiaa Ossama, Dr. Menna tu''Alllah Mohammad, Dr. Fatma Elhady, Dr. Abeer Ragab, Dr. Kerelos Zareef'
$ws.Range("G9").Value = 'Dr. Safa Hany, Dr. Shimaa Ashraf'
$ws.Range("G11").Value = 'Dr. Aya Saeed, Dr. Amal Awwad, Dr. Safa Hany'
$ws.Range("G12").Value = 'Dr. Eman M. Abo-Sakaya, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Dina Adel, Dr. Amira Ibrahim, Dr. Marina Youhanna'
$ws.Range("G13").Value = 'Dr. Amira Ibrahim, Dr. Esraa Mostafa, Dr. Yasmeena Fattoh'
$ws.Range("G15").Value = 'Dr. Rania Ahmad Youssef, Dr. Mohammad Safwat'
$ws.Range("G17").Value = 'Dr. Mohammad Safwat, Dr. Esraa Samy'
$ws.Range("G20").Value = 'Dr. Mariam Toma Gerges, Dr. Mohammad Safwat'
$ws.Range("G27").Value = 'Dr. Nourham Mostafa, Dr. Hana Amr'
$ws.Range("G28").Value = 'Dr. Aya Emad, Dr. Maryam Ashraf'
$ws.Range("G30").Value = 'Dr. Yassmen Ahmad, Dr. Aya Hanafy, Dr. Shorok Mohammad, Dr. Wafaa Ebida'
